$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 86, shifting existing rows 86-101 down to 87-102.
$ws.Rows.Item(86).Insert()

# Populate the new row 86 with the new record's values.
$ws.Cells.Item(86, 1).Value = 1
$ws.Cells.Item(86, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(86, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(86, 4).Value = 44641
$ws.Cells.Item(86, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(86, 5).Value = 15
$ws.Cells.Item(86, 6).Value = 100112008
$ws.Cells.Item(86, 7).Value = "Coliflor"
$ws.Cells.Item(86, 8).Value = "Sin especificar"
$ws.Cells.Item(86, 9).Value = "Segunda"
$ws.Cells.Item(86, 10).Value = 1300
$ws.Cells.Item(86, 11).Value = 550
$ws.Cells.Item(86, 12).Value = 600
$ws.Cells.Item(86, 13).Value = 575
$ws.Cells.Item(86, 14).Value = "$/unidad"
$ws.Cells.Item(86, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(86, 16).Value = 575
$ws.Cells.Item(86, 17).Value = 1
$ws.Cells.Item(86, 18).Value = "Hortaliza"
